$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (rows 2-51 => B:E), matching the
# "Updated cryptos list" GitHub Actions commit.
$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "23.124.04", "  -2.93%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.606.02", "  -2.87%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9996", "  -0.35%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9999", "  -0.15%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "301.98", "  -2.93%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.3785", "  -2.76%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3667", "  -3.91%  "),
    @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "49.92", "  -2.79%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.273", "  -5.66%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.08164", "  -3.69%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9996", "  -0.37%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "23.02", "  -4.09%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "6.630", "  -6.05%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001264", "  -3.79%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.417", "  -8.36%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.602.72", "  -3.08%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "92.11", "  -2.16%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06865", "  -2.02%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "18.33", "  -6.57%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.603", "  -5.80%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  -0.10%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "13.07", "  -4.82%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "23.126.38", "  -2.94%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.367", "  -2.81%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.774", "  -5.86%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "21.19", "  -3.99%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "150.02", "  -2.69%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "5.273", "  -3.04%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "133.14", "  -3.72%  "),
    @("WEMIXTOKEN", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.363", "  -5.23%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "6.969", "  -11.18%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.781.08", "  -3.03%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.9637", "  -4.65%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.07726", "  -5.62%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "6.342", "  -4.80%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02731", "  -6.06%  "),
    @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2556", "  -4.44%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "10.16", "  -5.74%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.08913", "  -2.74%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.372", "  -3.74%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.7118", "  -5.89%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "12.71", "  -5.97%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "15.30", "  -7.02%  "),
    @("Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.6642", "  -4.30%  "),
    @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "2.322", "  -5.19%  "),
    @("Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9991", "  +0.55%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "4.002", "  -2.63%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "132.10", "  -1.30%  "),
    @("Flow", "https://coinranking.com/coin/QQ0NCmjVq+flow-flow", "1.240", "  +1.06%  "),
    @("Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.07948", "  -4.12%  ")
)

# Column D holds price strings such as "23.124.04" / "0.9996" / "1.000"
# that must stay literal text (not be reinterpreted as numbers), so force
# the Text number format before writing, then clear the format back off
# so no stray style survives on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $row++
}

$ws.Range("D2:D51").ClearFormats()
